$wb = $excel.ActiveWorkbook

# --- Rename "Agency" sheet to "AgencyType" and update its header row ---
$wsAgency = $wb.Worksheets.Item("Agency")
$wsAgency.Name = "AgencyType"
$wsAgency.Range("A1").Value = "AgencyTypeID"
$wsAgency.Range("B1").Value = "AgencyTypeDescription"
$wsAgency.Columns.Item(1).ColumnWidth = 11.3
[void]$wsAgency.Range("I10").Select()

# --- Rename "Jurisdiction" sheet to "JurisdictionType" and update its header row ---
$wsJurisdiction = $wb.Worksheets.Item("Jurisdiction")
$wsJurisdiction.Name = "JurisdictionType"
$wsJurisdiction.Range("A1").Value = "JurisdictionTypeID"
$wsJurisdiction.Range("B1").Value = "JurisdictionTypeDescription"
[void]$wsJurisdiction.Range("D5").Select()

# --- Rename "EducationType" sheet to "EducationLevelType" and update its header row ---
$wsEducation = $wb.Worksheets.Item("EducationType")
$wsEducation.Name = "EducationLevelType"
$wsEducation.Range("A1").Value = "EducationLevelTypeID"
$wsEducation.Range("B1").Value = "EducationLevelTypeDescription"
$wsEducation.Columns.Item(2).ColumnWidth = 23.8
[void]$wsEducation.Range("E5").Select()
